# Fig/Intro/Intro_Sleep_Stages_PSD: fix the PSD (spectrogram) subplot's
# y-axis label on slide 1 — it was mistakenly labelled "Amplitude [µV]"
# (copy/pasted from the EEG amplitude axis above it) but the spectrogram's
# y-axis is actually frequency, so it should read "Frequency [Hz]".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "ZoneTexte 89" (shape #17 on this slide) is the vertical axis-label
# textbox that sits to the left of the bottom-left PSD/spectrogram plot.
$shp = $s.Shapes.Item("ZoneTexte 89")

# The textbox auto-fits its height to the text ("spAutoFit"). Remember the
# current (autofit) height so we can restore it after the text swap: the
# old and new labels are both 14 characters and wrap onto a single line,
# so the rendered height should not actually change.
$origHeight = $shp.Height

$shp.TextFrame.TextRange.Text = "Frequency [Hz]"

# Restore the height (nudged by a hair to survive the points<->EMU
# round-trip) so the shape's size stays byte-identical to the original.
$shp.Height = [single]($origHeight + 0.00002)
